# Auto-generated Excel COM-interop script applying the Maduin_Profits data refresh
# (scheduled runner update of currentAveragePrice / Leve price / profit columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1635.8636
$ws.Range("I12").Value = 1052.1052
$ws.Range("J12").Value = 5333
$ws.Range("K12").Value = 1052.1052
$ws.Range("L12").Value = 5333
$ws.Range("M12").Value = -882.1052
$ws.Range("N12").Value = -5673

$ws.Range("H15").Value = 584.9
$ws.Range("I15").Value = 584.9
$ws.Range("K15").Value = 1754.7
$ws.Range("M15").Value = -1585.7

$ws.Range("H18").Value = 2166
$ws.Range("J18").Value = 1000
$ws.Range("L18").Value = 1000
$ws.Range("N18").Value = -1568

$ws.Range("H62").Value = 5393.25
$ws.Range("I62").Value = 6335
$ws.Range("J62").Value = 4828.2
$ws.Range("K62").Value = 6335
$ws.Range("L62").Value = 4828.2
$ws.Range("M62").Value = -5711
$ws.Range("N62").Value = -6076.2

$ws.Range("H65").Value = 5393.25
$ws.Range("I65").Value = 6335
$ws.Range("J65").Value = 4828.2
$ws.Range("K65").Value = 31675
$ws.Range("L65").Value = 24141
$ws.Range("M65").Value = -28555
$ws.Range("N65").Value = -30381

$ws.Range("H88").Value = 1857.7727
$ws.Range("I88").Value = 832.6667
$ws.Range("J88").Value = 2019.6316
$ws.Range("K88").Value = 832.6667
$ws.Range("L88").Value = 2019.6316
$ws.Range("M88").Value = -426.6667
$ws.Range("N88").Value = -2831.6316

$ws.Range("H91").Value = 1857.7727
$ws.Range("I91").Value = 832.6667
$ws.Range("J91").Value = 2019.6316
$ws.Range("K91").Value = 832.6667
$ws.Range("L91").Value = 2019.6316
$ws.Range("M91").Value = 571.3333
$ws.Range("N91").Value = -4827.6316

$ws.Range("H93").Value = 29000
$ws.Range("J93").Value = 29000
$ws.Range("L93").Value = 29000
$ws.Range("N93").Value = -33992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3367.5386
$ws.Range("I45").Value = 1347.25
$ws.Range("K45").Value = 1347.25
$ws.Range("M45").Value = -970.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5551.375
$ws.Range("I105").Value = 4581
$ws.Range("J105").Value = 5874.8335
$ws.Range("K105").Value = 4581
$ws.Range("L105").Value = 5874.8335
$ws.Range("M105").Value = -2834
$ws.Range("N105").Value = -9368.833500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1073
$ws.Range("I16").Value = 1030.6666
$ws.Range("K16").Value = 1030.6666
$ws.Range("M16").Value = -743.6666

$ws.Range("H59").Value = 46665.668
$ws.Range("J59").Value = 46665.668
$ws.Range("L59").Value = 46665.668
$ws.Range("N59").Value = -48955.668

$ws.Range("H60").Value = 25746.75
$ws.Range("J60").Value = 25746.75
$ws.Range("L60").Value = 25746.75
$ws.Range("N60").Value = -26768.75

$ws.Range("H113").Value = 1073
$ws.Range("I113").Value = 1030.6666
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 1030.6666
$ws.Range("M113").Value = 1139.3334

$ws.Range("H122").Value = 962.2
$ws.Range("I122").Value = 1127.75
$ws.Range("J122").Value = 300
$ws.Range("K122").Value = 3383.25
$ws.Range("L122").Value = 900
$ws.Range("M122").Value = -933.25
$ws.Range("N122").Value = -5800

$ws.Range("H139").Value = 85890
$ws.Range("J139").Value = 85890
$ws.Range("L139").Value = 85890
$ws.Range("N139").Value = -96170

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 864
$ws.Range("I23").Value = 650
$ws.Range("J23").Value = 917.5
$ws.Range("K23").Value = 1950
$ws.Range("L23").Value = 2752.5
$ws.Range("M23").Value = -1715
$ws.Range("N23").Value = -3222.5

$ws.Range("H49").Value = 4300
$ws.Range("I49").Value = 5000
$ws.Range("J49").Value = 4125
$ws.Range("K49").Value = 15000
$ws.Range("L49").Value = 12375
$ws.Range("M49").Value = -14844
$ws.Range("N49").Value = -12687

$ws.Range("H103").Value = 2168.8
$ws.Range("I103").Value = 448
$ws.Range("J103").Value = 4750
$ws.Range("K103").Value = 1344
$ws.Range("L103").Value = 14250
$ws.Range("M103").Value = -465
$ws.Range("N103").Value = -16008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 4930
$ws.Range("I99").Value = 5917.5
$ws.Range("K99").Value = 5917.5
$ws.Range("M99").Value = -3671.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1579.7778
$ws.Range("J22").Value = 1530.4
$ws.Range("L22").Value = 1530.4
$ws.Range("N22").Value = -2120.4

$ws.Range("H26").Value = 5000
$ws.Range("J26").Value = 5000
$ws.Range("L26").Value = 5000
$ws.Range("N26").Value = -5590

$ws.Range("H27").Value = 1579.7778
$ws.Range("J27").Value = 1530.4
$ws.Range("L27").Value = 1530.4
$ws.Range("N27").Value = -1744.4

$ws.Range("H46").Value = 2161.762
$ws.Range("J46").Value = 2650
$ws.Range("L46").Value = 2650
$ws.Range("N46").Value = -3026

$ws.Range("H55").Value = 873.3043
$ws.Range("J55").Value = 1005.3333
$ws.Range("L55").Value = 1005.3333
$ws.Range("N55").Value = -1351.3333

$ws.Range("H100").Value = 2215.7144
$ws.Range("I100").Value = 1502.2
$ws.Range("K100").Value = 1502.2
$ws.Range("M100").Value = -961.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 39999
$ws.Range("J56").Value = 39999
$ws.Range("L56").Value = 39999
$ws.Range("N56").Value = -41427

$ws.Range("H80").Value = 26767.334
$ws.Range("J80").Value = 26767.334
$ws.Range("L80").Value = 26767.334
$ws.Range("N80").Value = -28763.334

$ws.Range("H83").Value = 26767.334
$ws.Range("J83").Value = 26767.334
$ws.Range("L83").Value = 80302.00199999999
$ws.Range("N83").Value = -90286.00199999999

$ws.Range("H122").Value = 1340.8572
$ws.Range("J122").Value = 847.5
$ws.Range("L122").Value = 2542.5
$ws.Range("N122").Value = -7442.5

$ws.Range("H132").Value = 84680.914
$ws.Range("I132").Value = 101317.1
$ws.Range("K132").Value = 303951.3
$ws.Range("M132").Value = -301421.3

$ws.Range("H136").Value = 2417.7856
$ws.Range("I136").Value = 1664.9
$ws.Range("J136").Value = 4300
$ws.Range("K136").Value = 4994.700000000001
$ws.Range("L136").Value = 12900
$ws.Range("M136").Value = -2444.700000000001
$ws.Range("N136").Value = -18000
